# Removed scenarios for Sorting Test as some scenarios would take too long to validate.
#
# On the "SortEmployeeTableColumns" sheet, remove the test-scenario rows for:
#   - Sort by Title on the Employee Page in Ascending Order      (row 6)
#   - Sort by Supervisor on the Employee Page in Descending Order (row 9)
#   - Sort by Location on the Employee Page in Ascending Order    (row 12)
#   - Sort by Vacation Days Left ... Descending Order              (row 15)
#   - Sort by Sick Days Left ... Descending Order                  (row 17)
#   - Sort by Floating Days Left ... Descending Order              (row 19)
#
# Deleting whole rows (bottom-most first so the remaining row numbers don't
# shift out from under us) lets Excel naturally re-pack xl/sheetData and
# drop the now-unused shared strings, which also renumbers the shared
# string indices referenced from the ShowEmployeesPerPage sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SortEmployeeTableColumns")

$ws.Rows.Item(19).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(6).Delete()

$ws.Activate()
$ws.Range("A6").Select() | Out-Null
